$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Duplicate row 2's uploaded record into row 3 ("data read") ---
$ws.Range("A3").Value = $ws.Range("A2").Value()
$ws.Range("B3").Value = $ws.Range("B2").Value()
$ws.Range("C3").Value = $ws.Range("C2").Value()
$ws.Hyperlinks.Add($ws.Range("C3"), "mailto:mithun.howlader222@gmail.com")

# --- Rows 2-10: newly uploaded / pending rows -> red, centered ---
$ws.Range("A2:B10").HorizontalAlignment = -4108
$ws.Range("A2:B10").Font.Color = 255

$ws.Range("C2:C10").Style = "Hyperlink"
$ws.Range("C2:C10").HorizontalAlignment = -4108
$ws.Range("C2:C10").Font.Color = 255

# --- Rows 11-21: pre-formatted blank rows for future uploads -> centered ---
$ws.Range("A11:B21").HorizontalAlignment = -4108

$ws.Range("C11:C21").Style = "Hyperlink"
$ws.Range("C11:C21").HorizontalAlignment = -4108

# --- Page setup / view bookkeeping ---
$ws.PageSetup.Orientation = 1
$ws.Range("G18").Select()
